# Commit: "Modified the function name Insert Cell value to Set Cell Value."
#
# Net effect on the workbook: two new sheets are inserted into the tab
# order (one blank, one carrying a "New Data" cell) and the pre-existing
# "New Sheet1" / "New Sheet" / "Introduction_Modified1" sheets end up
# shifted/renamed to make room for them, while "Introduction_Modified"
# and "Sheet" are untouched.
#
# Target tab order (name -> B1 content):
#   1. New Sheet3             -> "New Data"   (was "New Sheet1")
#   2. New Sheet2              -> (empty, new)
#   3. New Sheet1              -> "New Data"   (was "New Sheet")
#   4. New Sheet                -> "New Data"   (was "Introduction_Modified1")
#   5. Introduction_Modified1  -> "New Data"   (new)
#   6. Introduction_Modified   -> (empty, untouched)
#   7. Sheet                    -> (empty, untouched)

$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing sheets that shift position -------------
$sheet1 = $wb.Worksheets.Item("New Sheet1")
$sheet1.Name = "New Sheet3"

$sheet2 = $wb.Worksheets.Item("New Sheet")
$sheet2.Name = "New Sheet1"

$sheet3 = $wb.Worksheets.Item("Introduction_Modified1")
$sheet3.Name = "New Sheet"

# --- Step 2: insert a brand-new blank sheet right after "New Sheet3" ----
# (re-look-up the anchor sheet by name right before the call: sheet
# handles in this runtime track positional index, and stale handles can
# drift once a prior Add() shifts later tab positions)
$afterSheet3 = $wb.Worksheets.Item("New Sheet3")
$newSheet2 = $wb.Worksheets.Add($null, $afterSheet3)
$newSheet2.Name = "New Sheet2"

# --- Step 3: insert a brand-new sheet (with data) right after "New Sheet"
$afterSheetNew = $wb.Worksheets.Item("New Sheet")
$newIntroSheet = $wb.Worksheets.Add($null, $afterSheetNew)
$newIntroSheet.Name = "Introduction_Modified1"
$newIntroSheet.Range("B1").Value = "New Data"

# --- Step 4: restore the original active tab (position 1) ---------------
$wb.Worksheets.Item(1).Activate()
